$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows: SuperWhiteMinion (row 10) and SuperBlackMinion (row 11) ---
# Set the name cells first, then the class-path cells, so shared-string
# allocation order matches the target (names before paths).
$ws.Range("A10").Value = "SuperWhiteMinion"
$ws.Range("A11").Value = "SuperBlackMinion"
$ws.Range("B10").Value = "/Game/Character/Monster/SuperWhiteMinion/BSuperWhiteMinion_BP.BSuperWhiteMinion_BP_C"
$ws.Range("B11").Value = "/Game/Character/Monster/SuperBlackMinion/BSuperBlackMinion_BP.BSuperBlackMinion_BP_C"

$ws.Range("C10").Value = 100
$ws.Range("D10").Value = 70
$ws.Range("E10").Value = -100
$ws.Range("F10").Value = 30
$ws.Range("G10").Value = 100
$ws.Range("H10").Value = 200
$ws.Range("I10").Value = 30
$ws.Range("I10").NumberFormat = "0_);\(0\)"

$ws.Range("C11").Value = 100
$ws.Range("D11").Value = 70
$ws.Range("E11").Value = -100
$ws.Range("F11").Value = 50
$ws.Range("G11").Value = 200
$ws.Range("H11").Value = 200
$ws.Range("I11").Value = 50
$ws.Range("I11").NumberFormat = "0_);\(0\)"

# --- Row 7: "Rapter" typo fixed to "Raptor"; blueprint path corrected too ---
# B7 originally carries a quote-prefix ('text) style; a leading apostrophe
# keeps Excel marking the cell the same way instead of resetting its style.
$ws.Range("A7").Value = "Raptor"
$ws.Range("B7").Value = "'/Game/Character/Monster/Rapter/BRapter_BP.BRaptor_BP_C"

# --- Selection / camera move for the new boss rows ---
$ws.Range("A10").Select()

# --- Editor window geometry (best effort; mirrors the saved workbookView
# xWindow/yWindow/windowWidth/windowHeight from the authored change) ---
$win = $excel.ActiveWindow
$win.Left = 8490
$win.Top = 4095
$win.Width = 21600
$win.Height = 11385
